# Pharma_Society_Report.xlsx edit: rename sheet, update membership counts,
# and leave the final selection on D12 (matching the author's last-saved view).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab from "Sheet1" to "Report"
$ws.Name = "Report"

# Update the membership-count column (B) for each society
$ws.Range("B2").Value = 600   # FLASCO: 500 -> 600
$ws.Range("B3").Value = 500   # GASCO: 150 -> 500
$ws.Range("B5").Value = 176   # IOWA Oncology Society: 63 -> 176
$ws.Range("B6").Value = 400   # MOASC: 330 -> 400

# Match the saved cursor/selection position from the workbook being edited
$ws.Range("D12").Select() | Out-Null
